$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text, matching original inline string cells
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.931.01"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "3.528.08"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "603.61"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "143.42"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("D7").Value = "3.528.22"
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  -4.71%  "
$ws.Range("D12").Value = "0.405"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "4.135.83"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("E14").Value = "  -7.70%  "
$ws.Range("D15").Value = "28.27"
$ws.Range("E15").Value = "  -5.87%  "
$ws.Range("D16").Value = "3.520.52"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "65.823.40"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "10.91"
$ws.Range("E19").Value = "  -4.91%  "
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").Value = "14.53"
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("D22").Value = "419.81"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("D23").Value = "0.594"
$ws.Range("E23").Value = "  -4.20%  "
$ws.Range("D25").Value = "3.665.43"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -7.10%  "
$ws.Range("D28").Value = "2.46"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("D29").Value = "7.81"
$ws.Range("E29").Value = "  -5.24%  "
$ws.Range("D30").Value = "8.87"
$ws.Range("E30").Value = "  -4.90%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").Value = "3.533.67"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").Value = "24.15"
$ws.Range("E34").Value = "  -5.52%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -7.01%  "
$ws.Range("D37").Value = "7.50"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").Value = "175.54"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("E39").Value = "  -6.93%  "
$ws.Range("D40").Value = "5.22"
$ws.Range("E40").Value = "  -7.32%  "
$ws.Range("D41").Value = "0.0813"
$ws.Range("E41").Value = "  -5.20%  "
$ws.Range("D42").Value = "0.856"
$ws.Range("E42").Value = "  -4.65%  "
$ws.Range("D43").Value = "4.94"
$ws.Range("E43").Value = "  -5.62%  "
$ws.Range("D44").Value = "45.42"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").Value = "1.77"
$ws.Range("E45").Value = "  -7.49%  "
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -9.13%  "
$ws.Range("D48").Value = "23.13"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("D51").Value = "0.902"
$ws.Range("E51").Value = "  -5.29%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "7.01"
$ws.Range("E49").Value = "  -2.93%  "

$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "1.11"
$ws.Range("E50").Value = "  -6.93%  "
